$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C width (descr) ---
$ws.Columns.Item(3).ColumnWidth = 65.8

# --- Phase 1: seed new shared strings in the exact order they first appear ---
# (so sharedStrings.xml ends up with the same index assignment as the target workbook)
$ws.Cells.Item(2495, "B").Value = 'ROOT_CERT_EXPIRY_TEMPLATE'
$ws.Cells.Item(2495, "C").Value = 'Template for root certificate expiry'
$ws.Cells.Item(2496, "C").Value = 'Modèle d''expiration du certificat racine'
$ws.Cells.Item(2497, "C").Value = 'نموذج لانتهاء صلاحية شهادة الجذر'
$ws.Cells.Item(2501, "B").Value = 'INTERMEDIATE_CERT_EXPIRY_TEMPLATE'
$ws.Cells.Item(2501, "C").Value = 'Template for intermediate certificate expiry'
$ws.Cells.Item(2502, "C").Value = 'Modèle d''expiration de certificat intermédiaire'
$ws.Cells.Item(2503, "C").Value = 'نموذج انتهاء صلاحية الشهادة المتوسطة'
$ws.Cells.Item(2507, "B").Value = 'PARTNER_CERT_EXPIRY_TEMPLATE'
$ws.Cells.Item(2507, "C").Value = 'Template for partner certificate expiry'
$ws.Cells.Item(2508, "C").Value = 'Modèle d''expiration du certificat de partenaire'
$ws.Cells.Item(2509, "C").Value = 'نموذج انتهاء صلاحية شهادة الشريك'
$ws.Cells.Item(2513, "B").Value = 'WEEKLY_SUMMARY_TEMPLATE'
$ws.Cells.Item(2513, "C").Value = 'Template for weekly summary notifications'
$ws.Cells.Item(2514, "C").Value = 'Modèle pour les notifications récapitulatives hebdomadaires'
$ws.Cells.Item(2515, "C").Value = 'نموذج لإشعارات الملخص الأسبوعية'
$ws.Cells.Item(2519, "B").Value = 'ROOT_CERT_EXPIRY_SUBJECT'
$ws.Cells.Item(2519, "C").Value = 'Subject template for root certificate expiry'
$ws.Cells.Item(2520, "C").Value = 'Modèle de sujet pour l''expiration du certificat racine'
$ws.Cells.Item(2521, "C").Value = 'نموذج موضوعي لانتهاء صلاحية شهادة الجذر'
$ws.Cells.Item(2525, "B").Value = 'INTERMEDIATE_CERT_EXPIRY_SUBJECT'
$ws.Cells.Item(2525, "C").Value = 'Subject template for intermediate certificate expiry'
$ws.Cells.Item(2526, "C").Value = 'Modèle de sujet pour l''expiration du certificat intermédiaire'
$ws.Cells.Item(2527, "C").Value = 'نموذج موضوعي لانتهاء صلاحية الشهادة المتوسطة'
$ws.Cells.Item(2531, "B").Value = 'PARTNER_CERT_EXPIRY_SUBJECT'
$ws.Cells.Item(2531, "C").Value = 'Subject template for partner certificate expiry'
$ws.Cells.Item(2532, "C").Value = 'Modèle de sujet pour l''expiration du certificat du partenaire'
$ws.Cells.Item(2533, "C").Value = 'نموذج موضوعي لانتهاء صلاحية شهادة الشريك'
$ws.Cells.Item(2537, "B").Value = 'WEEKLY_SUMMARY_SUBJECT'
$ws.Cells.Item(2537, "C").Value = 'Subject template for weekly summary notifications'
$ws.Cells.Item(2538, "C").Value = 'Modèle de sujet pour les notifications récapitulatives hebdomadaires'
$ws.Cells.Item(2539, "C").Value = 'قالب موضوعي لإشعارات الملخص الأسبوعي'
$ws.Cells.Item(2500, "C").Value = 'ரூட் சான்றிதழ் காலாவதிக்கான டெம்ப்ளேட்'
$ws.Cells.Item(2498, "C").Value = 'रूट प्रमाणपत्र समाप्ति के लिए टेम्पलेट'
$ws.Cells.Item(2499, "C").Value = 'ಮೂಲ ಪ್ರಮಾಣಪತ್ರದ ಮುಕ್ತಾಯ ದಿನಾಂಕದ ಟೆಂಪ್ಲೇಟ್'
$ws.Cells.Item(2505, "C").Value = 'ಮಧ್ಯಂತರ ಪ್ರಮಾಣಪತ್ರ ಮುಕ್ತಾಯಕ್ಕಾಗಿ ಟೆಂಪ್ಲೇಟ್'
$ws.Cells.Item(2504, "C").Value = 'मध्यवर्ती प्रमाणपत्र समाप्ति के लिए टेम्पलेट'
$ws.Cells.Item(2506, "C").Value = 'இடைநிலை சான்றிதழ் காலாவதிக்கான டெம்ப்ளேட்'
$ws.Cells.Item(2512, "C").Value = 'கூட்டாளர் சான்றிதழ் காலாவதிக்கான டெம்ப்ளேட்'
$ws.Cells.Item(2511, "C").Value = 'ಪಾಲುದಾರ ಪ್ರಮಾಣಪತ್ರದ ಮುಕ್ತಾಯ ದಿನಾಂಕದ ಟೆಂಪ್ಲೇಟ್'
$ws.Cells.Item(2510, "C").Value = 'भागीदार प्रमाणपत्र समाप्ति के लिए टेम्पलेट'
$ws.Cells.Item(2516, "C").Value = 'साप्ताहिक सारांश अधिसूचनाओं के लिए टेम्पलेट'
$ws.Cells.Item(2517, "C").Value = 'ವಾರದ ಸಾರಾಂಶ ಅಧಿಸೂಚನೆಗಳಿಗಾಗಿ ಟೆಂಪ್ಲೇಟ್'
$ws.Cells.Item(2518, "C").Value = 'வாராந்திர சுருக்க அறிவிப்புகளுக்கான டெம்ப்ளேட்'
$ws.Cells.Item(2524, "C").Value = 'மூலச் சான்றிதழ் காலாவதிக்கான பொருள் டெம்ப்ளேட்'
$ws.Cells.Item(2522, "C").Value = 'रूट प्रमाणपत्र समाप्ति के लिए विषय टेम्पलेट'
$ws.Cells.Item(2523, "C").Value = 'ಮೂಲ ಪ್ರಮಾಣಪತ್ರದ ಮುಕ್ತಾಯ ದಿನಾಂಕದ ವಿಷಯ ಟೆಂಪ್ಲೇಟ್'
$ws.Cells.Item(2529, "C").Value = 'ಮಧ್ಯಂತರ ಪ್ರಮಾಣಪತ್ರ ಮುಕ್ತಾಯಕ್ಕಾಗಿ ವಿಷಯ ಟೆಂಪ್ಲೇಟ್'
$ws.Cells.Item(2528, "C").Value = 'इंटरमीडिएट प्रमाणपत्र समाप्ति के लिए विषय टेम्पलेट'
$ws.Cells.Item(2530, "C").Value = 'இடைநிலை சான்றிதழ் காலாவதிக்கான பொருள் வார்ப்புரு'
$ws.Cells.Item(2536, "C").Value = 'கூட்டாளர் சான்றிதழ் காலாவதிக்கான பொருள் டெம்ப்ளேட்'
$ws.Cells.Item(2535, "C").Value = 'ಪಾಲುದಾರ ಪ್ರಮಾಣಪತ್ರದ ಮುಕ್ತಾಯದ ವಿಷಯ ಟೆಂಪ್ಲೇಟ್'
$ws.Cells.Item(2534, "C").Value = 'भागीदार प्रमाणपत्र समाप्ति के लिए विषय टेम्पलेट'
$ws.Cells.Item(2540, "C").Value = 'साप्ताहिक सारांश अधिसूचनाओं के लिए विषय टेम्पलेट'
$ws.Cells.Item(2541, "C").Value = 'ವಾರದ ಸಾರಾಂಶ ಅಧಿಸೂಚನೆಗಳಿಗಾಗಿ ವಿಷಯ ಟೆಂಪ್ಲೇಟ್'
$ws.Cells.Item(2542, "C").Value = 'வாராந்திர சுருக்க அறிவிப்புகளுக்கான தலைப்பு டெம்ப்ளேட்'

# --- Phase 2: fill in the remaining cells for the new template rows (2495-2542) ---
# row 2495
$ws.Cells.Item(2495, "A").Value = 'eng'
$ws.Cells.Item(2495, "B").Value = 'ROOT_CERT_EXPIRY_TEMPLATE'
$ws.Cells.Item(2495, "C").Value = 'Template for root certificate expiry'
$ws.Cells.Item(2, 4).Copy($ws.Cells.Item(2495, 4)) | Out-Null
$ws.Cells.Item(2495, "E").Value = 'admin'
$ws.Cells.Item(2495, "F").Value = 'now()'

# row 2496
$ws.Cells.Item(2496, "A").Value = 'fra'
$ws.Cells.Item(2496, "B").Value = 'ROOT_CERT_EXPIRY_TEMPLATE'
$ws.Cells.Item(2496, "C").Value = 'Modèle d''expiration du certificat racine'
$ws.Cells.Item(2, 4).Copy($ws.Cells.Item(2496, 4)) | Out-Null
$ws.Cells.Item(2496, "E").Value = 'admin'
$ws.Cells.Item(2496, "F").Value = 'now()'

# row 2497
$ws.Cells.Item(2497, "A").Value = 'ara'
$ws.Cells.Item(2497, "B").Value = 'ROOT_CERT_EXPIRY_TEMPLATE'
$ws.Cells.Item(2497, "C").Value = 'نموذج لانتهاء صلاحية شهادة الجذر'
$ws.Cells.Item(2, 4).Copy($ws.Cells.Item(2497, 4)) | Out-Null
$ws.Cells.Item(2497, "E").Value = 'admin'
$ws.Cells.Item(2497, "F").Value = 'now()'

# row 2498
$ws.Cells.Item(2498, "A").Value = 'hin'
$ws.Cells.Item(2498, "B").Value = 'ROOT_CERT_EXPIRY_TEMPLATE'
$ws.Cells.Item(2498, "C").Value = 'रूट प्रमाणपत्र समाप्ति के लिए टेम्पलेट'
$ws.Cells.Item(2, 4).Copy($ws.Cells.Item(2498, 4)) | Out-Null
$ws.Cells.Item(2498, "E").Value = 'admin'
$ws.Cells.Item(2498, "F").Value = 'now()'

# row 2499
$ws.Cells.Item(2499, "A").Value = 'kan'
$ws.Cells.Item(2499, "B").Value = 'ROOT_CERT_EXPIRY_TEMPLATE'
$ws.Cells.Item(2499, "C").Value = 'ಮೂಲ ಪ್ರಮಾಣಪತ್ರದ ಮುಕ್ತಾಯ ದಿನಾಂಕದ ಟೆಂಪ್ಲೇಟ್'
$ws.Cells.Item(2, 4).Copy($ws.Cells.Item(2499, 4)) | Out-Null
$ws.Cells.Item(2499, "E").Value = 'admin'
$ws.Cells.Item(2499, "F").Value = 'now()'

# row 2500
$ws.Cells.Item(2500, "A").Value = 'tam'
$ws.Cells.Item(2500, "B").Value = 'ROOT_CERT_EXPIRY_TEMPLATE'
$ws.Cells.Item(2500, "C").Value = 'ரூட் சான்றிதழ் காலாவதிக்கான டெம்ப்ளேட்'
$ws.Cells.Item(2, 4).Copy($ws.Cells.Item(2500, 4)) | Out-Null
$ws.Cells.Item(2500, "E").Value = 'admin'
$ws.Cells.Item(2500, "F").Value = 'now()'

# row 2501
$ws.Cells.Item(2501, "A").Value = 'eng'
$ws.Cells.Item(2501, "B").Value = 'INTERMEDIATE_CERT_EXPIRY_TEMPLATE'
$ws.Cells.Item(2501, "C").Value = 'Template for intermediate certificate expiry'
$ws.Cells.Item(2, 4).Copy($ws.Cells.Item(2501, 4)) | Out-Null
$ws.Cells.Item(2501, "E").Value = 'admin'
$ws.Cells.Item(2501, "F").Value = 'now()'

# row 2502
$ws.Cells.Item(2502, "A").Value = 'fra'
$ws.Cells.Item(2502, "B").Value = 'INTERMEDIATE_CERT_EXPIRY_TEMPLATE'
$ws.Cells.Item(2502, "C").Value = 'Modèle d''expiration de certificat intermédiaire'
$ws.Cells.Item(2, 4).Copy($ws.Cells.Item(2502, 4)) | Out-Null
$ws.Cells.Item(2502, "E").Value = 'admin'
$ws.Cells.Item(2502, "F").Value = 'now()'

# row 2503
$ws.Cells.Item(2503, "A").Value = 'ara'
$ws.Cells.Item(2503, "B").Value = 'INTERMEDIATE_CERT_EXPIRY_TEMPLATE'
$ws.Cells.Item(2503, "C").Value = 'نموذج انتهاء صلاحية الشهادة المتوسطة'
$ws.Cells.Item(2, 4).Copy($ws.Cells.Item(2503, 4)) | Out-Null
$ws.Cells.Item(2503, "E").Value = 'admin'
$ws.Cells.Item(2503, "F").Value = 'now()'

# row 2504
$ws.Cells.Item(2504, "A").Value = 'hin'
$ws.Cells.Item(2504, "B").Value = 'INTERMEDIATE_CERT_EXPIRY_TEMPLATE'
$ws.Cells.Item(2504, "C").Value = 'मध्यवर्ती प्रमाणपत्र समाप्ति के लिए टेम्पलेट'
$ws.Cells.Item(2, 4).Copy($ws.Cells.Item(2504, 4)) | Out-Null
$ws.Cells.Item(2504, "E").Value = 'admin'
$ws.Cells.Item(2504, "F").Value = 'now()'

# row 2505
$ws.Cells.Item(2505, "A").Value = 'kan'
$ws.Cells.Item(2505, "B").Value = 'INTERMEDIATE_CERT_EXPIRY_TEMPLATE'
$ws.Cells.Item(2505, "C").Value = 'ಮಧ್ಯಂತರ ಪ್ರಮಾಣಪತ್ರ ಮುಕ್ತಾಯಕ್ಕಾಗಿ ಟೆಂಪ್ಲೇಟ್'
$ws.Cells.Item(2, 4).Copy($ws.Cells.Item(2505, 4)) | Out-Null
$ws.Cells.Item(2505, "E").Value = 'admin'
$ws.Cells.Item(2505, "F").Value = 'now()'

# row 2506
$ws.Cells.Item(2506, "A").Value = 'tam'
$ws.Cells.Item(2506, "B").Value = 'INTERMEDIATE_CERT_EXPIRY_TEMPLATE'
$ws.Cells.Item(2506, "C").Value = 'இடைநிலை சான்றிதழ் காலாவதிக்கான டெம்ப்ளேட்'
$ws.Cells.Item(2, 4).Copy($ws.Cells.Item(2506, 4)) | Out-Null
$ws.Cells.Item(2506, "E").Value = 'admin'
$ws.Cells.Item(2506, "F").Value = 'now()'

# row 2507
$ws.Cells.Item(2507, "A").Value = 'eng'
$ws.Cells.Item(2507, "B").Value = 'PARTNER_CERT_EXPIRY_TEMPLATE'
$ws.Cells.Item(2507, "C").Value = 'Template for partner certificate expiry'
$ws.Cells.Item(2, 4).Copy($ws.Cells.Item(2507, 4)) | Out-Null
$ws.Cells.Item(2507, "E").Value = 'admin'
$ws.Cells.Item(2507, "F").Value = 'now()'

# row 2508
$ws.Cells.Item(2508, "A").Value = 'fra'
$ws.Cells.Item(2508, "B").Value = 'PARTNER_CERT_EXPIRY_TEMPLATE'
$ws.Cells.Item(2508, "C").Value = 'Modèle d''expiration du certificat de partenaire'
$ws.Cells.Item(2, 4).Copy($ws.Cells.Item(2508, 4)) | Out-Null
$ws.Cells.Item(2508, "E").Value = 'admin'
$ws.Cells.Item(2508, "F").Value = 'now()'

# row 2509
$ws.Cells.Item(2509, "A").Value = 'ara'
$ws.Cells.Item(2509, "B").Value = 'PARTNER_CERT_EXPIRY_TEMPLATE'
$ws.Cells.Item(2509, "C").Value = 'نموذج انتهاء صلاحية شهادة الشريك'
$ws.Cells.Item(2, 4).Copy($ws.Cells.Item(2509, 4)) | Out-Null
$ws.Cells.Item(2509, "E").Value = 'admin'
$ws.Cells.Item(2509, "F").Value = 'now()'

# row 2510
$ws.Cells.Item(2510, "A").Value = 'hin'
$ws.Cells.Item(2510, "B").Value = 'PARTNER_CERT_EXPIRY_TEMPLATE'
$ws.Cells.Item(2510, "C").Value = 'भागीदार प्रमाणपत्र समाप्ति के लिए टेम्पलेट'
$ws.Cells.Item(2, 4).Copy($ws.Cells.Item(2510, 4)) | Out-Null
$ws.Cells.Item(2510, "E").Value = 'admin'
$ws.Cells.Item(2510, "F").Value = 'now()'

# row 2511
$ws.Cells.Item(2511, "A").Value = 'kan'
$ws.Cells.Item(2511, "B").Value = 'PARTNER_CERT_EXPIRY_TEMPLATE'
$ws.Cells.Item(2511, "C").Value = 'ಪಾಲುದಾರ ಪ್ರಮಾಣಪತ್ರದ ಮುಕ್ತಾಯ ದಿನಾಂಕದ ಟೆಂಪ್ಲೇಟ್'
$ws.Cells.Item(2, 4).Copy($ws.Cells.Item(2511, 4)) | Out-Null
$ws.Cells.Item(2511, "E").Value = 'admin'
$ws.Cells.Item(2511, "F").Value = 'now()'

# row 2512
$ws.Cells.Item(2512, "A").Value = 'tam'
$ws.Cells.Item(2512, "B").Value = 'PARTNER_CERT_EXPIRY_TEMPLATE'
$ws.Cells.Item(2512, "C").Value = 'கூட்டாளர் சான்றிதழ் காலாவதிக்கான டெம்ப்ளேட்'
$ws.Cells.Item(2, 4).Copy($ws.Cells.Item(2512, 4)) | Out-Null
$ws.Cells.Item(2512, "E").Value = 'admin'
$ws.Cells.Item(2512, "F").Value = 'now()'

# row 2513
$ws.Cells.Item(2513, "A").Value = 'eng'
$ws.Cells.Item(2513, "B").Value = 'WEEKLY_SUMMARY_TEMPLATE'
$ws.Cells.Item(2513, "C").Value = 'Template for weekly summary notifications'
$ws.Cells.Item(2, 4).Copy($ws.Cells.Item(2513, 4)) | Out-Null
$ws.Cells.Item(2513, "E").Value = 'admin'
$ws.Cells.Item(2513, "F").Value = 'now()'

# row 2514
$ws.Cells.Item(2514, "A").Value = 'fra'
$ws.Cells.Item(2514, "B").Value = 'WEEKLY_SUMMARY_TEMPLATE'
$ws.Cells.Item(2514, "C").Value = 'Modèle pour les notifications récapitulatives hebdomadaires'
$ws.Cells.Item(2, 4).Copy($ws.Cells.Item(2514, 4)) | Out-Null
$ws.Cells.Item(2514, "E").Value = 'admin'
$ws.Cells.Item(2514, "F").Value = 'now()'

# row 2515
$ws.Cells.Item(2515, "A").Value = 'ara'
$ws.Cells.Item(2515, "B").Value = 'WEEKLY_SUMMARY_TEMPLATE'
$ws.Cells.Item(2515, "C").Value = 'نموذج لإشعارات الملخص الأسبوعية'
$ws.Cells.Item(2, 4).Copy($ws.Cells.Item(2515, 4)) | Out-Null
$ws.Cells.Item(2515, "E").Value = 'admin'
$ws.Cells.Item(2515, "F").Value = 'now()'

# row 2516
$ws.Cells.Item(2516, "A").Value = 'hin'
$ws.Cells.Item(2516, "B").Value = 'WEEKLY_SUMMARY_TEMPLATE'
$ws.Cells.Item(2516, "C").Value = 'साप्ताहिक सारांश अधिसूचनाओं के लिए टेम्पलेट'
$ws.Cells.Item(2, 4).Copy($ws.Cells.Item(2516, 4)) | Out-Null
$ws.Cells.Item(2516, "E").Value = 'admin'
$ws.Cells.Item(2516, "F").Value = 'now()'

# row 2517
$ws.Cells.Item(2517, "A").Value = 'kan'
$ws.Cells.Item(2517, "B").Value = 'WEEKLY_SUMMARY_TEMPLATE'
$ws.Cells.Item(2517, "C").Value = 'ವಾರದ ಸಾರಾಂಶ ಅಧಿಸೂಚನೆಗಳಿಗಾಗಿ ಟೆಂಪ್ಲೇಟ್'
$ws.Cells.Item(2, 4).Copy($ws.Cells.Item(2517, 4)) | Out-Null
$ws.Cells.Item(2517, "E").Value = 'admin'
$ws.Cells.Item(2517, "F").Value = 'now()'

# row 2518
$ws.Cells.Item(2518, "A").Value = 'tam'
$ws.Cells.Item(2518, "B").Value = 'WEEKLY_SUMMARY_TEMPLATE'
$ws.Cells.Item(2518, "C").Value = 'வாராந்திர சுருக்க அறிவிப்புகளுக்கான டெம்ப்ளேட்'
$ws.Cells.Item(2, 4).Copy($ws.Cells.Item(2518, 4)) | Out-Null
$ws.Cells.Item(2518, "E").Value = 'admin'
$ws.Cells.Item(2518, "F").Value = 'now()'

# row 2519
$ws.Cells.Item(2519, "A").Value = 'eng'
$ws.Cells.Item(2519, "B").Value = 'ROOT_CERT_EXPIRY_SUBJECT'
$ws.Cells.Item(2519, "C").Value = 'Subject template for root certificate expiry'
$ws.Cells.Item(2, 4).Copy($ws.Cells.Item(2519, 4)) | Out-Null
$ws.Cells.Item(2519, "E").Value = 'admin'
$ws.Cells.Item(2519, "F").Value = 'now()'

# row 2520
$ws.Cells.Item(2520, "A").Value = 'fra'
$ws.Cells.Item(2520, "B").Value = 'ROOT_CERT_EXPIRY_SUBJECT'
$ws.Cells.Item(2520, "C").Value = 'Modèle de sujet pour l''expiration du certificat racine'
$ws.Cells.Item(2, 4).Copy($ws.Cells.Item(2520, 4)) | Out-Null
$ws.Cells.Item(2520, "E").Value = 'admin'
$ws.Cells.Item(2520, "F").Value = 'now()'

# row 2521
$ws.Cells.Item(2521, "A").Value = 'ara'
$ws.Cells.Item(2521, "B").Value = 'ROOT_CERT_EXPIRY_SUBJECT'
$ws.Cells.Item(2521, "C").Value = 'نموذج موضوعي لانتهاء صلاحية شهادة الجذر'
$ws.Cells.Item(2, 4).Copy($ws.Cells.Item(2521, 4)) | Out-Null
$ws.Cells.Item(2521, "E").Value = 'admin'
$ws.Cells.Item(2521, "F").Value = 'now()'

# row 2522
$ws.Cells.Item(2522, "A").Value = 'hin'
$ws.Cells.Item(2522, "B").Value = 'ROOT_CERT_EXPIRY_SUBJECT'
$ws.Cells.Item(2522, "C").Value = 'रूट प्रमाणपत्र समाप्ति के लिए विषय टेम्पलेट'
$ws.Cells.Item(2, 4).Copy($ws.Cells.Item(2522, 4)) | Out-Null
$ws.Cells.Item(2522, "E").Value = 'admin'
$ws.Cells.Item(2522, "F").Value = 'now()'

# row 2523
$ws.Cells.Item(2523, "A").Value = 'kan'
$ws.Cells.Item(2523, "B").Value = 'ROOT_CERT_EXPIRY_SUBJECT'
$ws.Cells.Item(2523, "C").Value = 'ಮೂಲ ಪ್ರಮಾಣಪತ್ರದ ಮುಕ್ತಾಯ ದಿನಾಂಕದ ವಿಷಯ ಟೆಂಪ್ಲೇಟ್'
$ws.Cells.Item(2, 4).Copy($ws.Cells.Item(2523, 4)) | Out-Null
$ws.Cells.Item(2523, "E").Value = 'admin'
$ws.Cells.Item(2523, "F").Value = 'now()'

# row 2524
$ws.Cells.Item(2524, "A").Value = 'tam'
$ws.Cells.Item(2524, "B").Value = 'ROOT_CERT_EXPIRY_SUBJECT'
$ws.Cells.Item(2524, "C").Value = 'மூலச் சான்றிதழ் காலாவதிக்கான பொருள் டெம்ப்ளேட்'
$ws.Cells.Item(2, 4).Copy($ws.Cells.Item(2524, 4)) | Out-Null
$ws.Cells.Item(2524, "E").Value = 'admin'
$ws.Cells.Item(2524, "F").Value = 'now()'

# row 2525
$ws.Cells.Item(2525, "A").Value = 'eng'
$ws.Cells.Item(2525, "B").Value = 'INTERMEDIATE_CERT_EXPIRY_SUBJECT'
$ws.Cells.Item(2525, "C").Value = 'Subject template for intermediate certificate expiry'
$ws.Cells.Item(2, 4).Copy($ws.Cells.Item(2525, 4)) | Out-Null
$ws.Cells.Item(2525, "E").Value = 'admin'
$ws.Cells.Item(2525, "F").Value = 'now()'

# row 2526
$ws.Cells.Item(2526, "A").Value = 'fra'
$ws.Cells.Item(2526, "B").Value = 'INTERMEDIATE_CERT_EXPIRY_SUBJECT'
$ws.Cells.Item(2526, "C").Value = 'Modèle de sujet pour l''expiration du certificat intermédiaire'
$ws.Cells.Item(2, 4).Copy($ws.Cells.Item(2526, 4)) | Out-Null
$ws.Cells.Item(2526, "E").Value = 'admin'
$ws.Cells.Item(2526, "F").Value = 'now()'

# row 2527
$ws.Cells.Item(2527, "A").Value = 'ara'
$ws.Cells.Item(2527, "B").Value = 'INTERMEDIATE_CERT_EXPIRY_SUBJECT'
$ws.Cells.Item(2527, "C").Value = 'نموذج موضوعي لانتهاء صلاحية الشهادة المتوسطة'
$ws.Cells.Item(2, 4).Copy($ws.Cells.Item(2527, 4)) | Out-Null
$ws.Cells.Item(2527, "E").Value = 'admin'
$ws.Cells.Item(2527, "F").Value = 'now()'

# row 2528
$ws.Cells.Item(2528, "A").Value = 'hin'
$ws.Cells.Item(2528, "B").Value = 'INTERMEDIATE_CERT_EXPIRY_SUBJECT'
$ws.Cells.Item(2528, "C").Value = 'इंटरमीडिएट प्रमाणपत्र समाप्ति के लिए विषय टेम्पलेट'
$ws.Cells.Item(2, 4).Copy($ws.Cells.Item(2528, 4)) | Out-Null
$ws.Cells.Item(2528, "E").Value = 'admin'
$ws.Cells.Item(2528, "F").Value = 'now()'

# row 2529
$ws.Cells.Item(2529, "A").Value = 'kan'
$ws.Cells.Item(2529, "B").Value = 'INTERMEDIATE_CERT_EXPIRY_SUBJECT'
$ws.Cells.Item(2529, "C").Value = 'ಮಧ್ಯಂತರ ಪ್ರಮಾಣಪತ್ರ ಮುಕ್ತಾಯಕ್ಕಾಗಿ ವಿಷಯ ಟೆಂಪ್ಲೇಟ್'
$ws.Cells.Item(2, 4).Copy($ws.Cells.Item(2529, 4)) | Out-Null
$ws.Cells.Item(2529, "E").Value = 'admin'
$ws.Cells.Item(2529, "F").Value = 'now()'

# row 2530
$ws.Cells.Item(2530, "A").Value = 'tam'
$ws.Cells.Item(2530, "B").Value = 'INTERMEDIATE_CERT_EXPIRY_SUBJECT'
$ws.Cells.Item(2530, "C").Value = 'இடைநிலை சான்றிதழ் காலாவதிக்கான பொருள் வார்ப்புரு'
$ws.Cells.Item(2, 4).Copy($ws.Cells.Item(2530, 4)) | Out-Null
$ws.Cells.Item(2530, "E").Value = 'admin'
$ws.Cells.Item(2530, "F").Value = 'now()'

# row 2531
$ws.Cells.Item(2531, "A").Value = 'eng'
$ws.Cells.Item(2531, "B").Value = 'PARTNER_CERT_EXPIRY_SUBJECT'
$ws.Cells.Item(2531, "C").Value = 'Subject template for partner certificate expiry'
$ws.Cells.Item(2, 4).Copy($ws.Cells.Item(2531, 4)) | Out-Null
$ws.Cells.Item(2531, "E").Value = 'admin'
$ws.Cells.Item(2531, "F").Value = 'now()'

# row 2532
$ws.Cells.Item(2532, "A").Value = 'fra'
$ws.Cells.Item(2532, "B").Value = 'PARTNER_CERT_EXPIRY_SUBJECT'
$ws.Cells.Item(2532, "C").Value = 'Modèle de sujet pour l''expiration du certificat du partenaire'
$ws.Cells.Item(2, 4).Copy($ws.Cells.Item(2532, 4)) | Out-Null
$ws.Cells.Item(2532, "E").Value = 'admin'
$ws.Cells.Item(2532, "F").Value = 'now()'

# row 2533
$ws.Cells.Item(2533, "A").Value = 'ara'
$ws.Cells.Item(2533, "B").Value = 'PARTNER_CERT_EXPIRY_SUBJECT'
$ws.Cells.Item(2533, "C").Value = 'نموذج موضوعي لانتهاء صلاحية شهادة الشريك'
$ws.Cells.Item(2, 4).Copy($ws.Cells.Item(2533, 4)) | Out-Null
$ws.Cells.Item(2533, "E").Value = 'admin'
$ws.Cells.Item(2533, "F").Value = 'now()'

# row 2534
$ws.Cells.Item(2534, "A").Value = 'hin'
$ws.Cells.Item(2534, "B").Value = 'PARTNER_CERT_EXPIRY_SUBJECT'
$ws.Cells.Item(2534, "C").Value = 'भागीदार प्रमाणपत्र समाप्ति के लिए विषय टेम्पलेट'
$ws.Cells.Item(2, 4).Copy($ws.Cells.Item(2534, 4)) | Out-Null
$ws.Cells.Item(2534, "E").Value = 'admin'
$ws.Cells.Item(2534, "F").Value = 'now()'

# row 2535
$ws.Cells.Item(2535, "A").Value = 'kan'
$ws.Cells.Item(2535, "B").Value = 'PARTNER_CERT_EXPIRY_SUBJECT'
$ws.Cells.Item(2535, "C").Value = 'ಪಾಲುದಾರ ಪ್ರಮಾಣಪತ್ರದ ಮುಕ್ತಾಯದ ವಿಷಯ ಟೆಂಪ್ಲೇಟ್'
$ws.Cells.Item(2, 4).Copy($ws.Cells.Item(2535, 4)) | Out-Null
$ws.Cells.Item(2535, "E").Value = 'admin'
$ws.Cells.Item(2535, "F").Value = 'now()'

# row 2536
$ws.Cells.Item(2536, "A").Value = 'tam'
$ws.Cells.Item(2536, "B").Value = 'PARTNER_CERT_EXPIRY_SUBJECT'
$ws.Cells.Item(2536, "C").Value = 'கூட்டாளர் சான்றிதழ் காலாவதிக்கான பொருள் டெம்ப்ளேட்'
$ws.Cells.Item(2, 4).Copy($ws.Cells.Item(2536, 4)) | Out-Null
$ws.Cells.Item(2536, "E").Value = 'admin'
$ws.Cells.Item(2536, "F").Value = 'now()'

# row 2537
$ws.Cells.Item(2537, "A").Value = 'eng'
$ws.Cells.Item(2537, "B").Value = 'WEEKLY_SUMMARY_SUBJECT'
$ws.Cells.Item(2537, "C").Value = 'Subject template for weekly summary notifications'
$ws.Cells.Item(2, 4).Copy($ws.Cells.Item(2537, 4)) | Out-Null
$ws.Cells.Item(2537, "E").Value = 'admin'
$ws.Cells.Item(2537, "F").Value = 'now()'
$ws.Rows.Item(2537).RowHeight = 14.25

# row 2538
$ws.Cells.Item(2538, "A").Value = 'fra'
$ws.Cells.Item(2538, "B").Value = 'WEEKLY_SUMMARY_SUBJECT'
$ws.Cells.Item(2538, "C").Value = 'Modèle de sujet pour les notifications récapitulatives hebdomadaires'
$ws.Cells.Item(2, 4).Copy($ws.Cells.Item(2538, 4)) | Out-Null
$ws.Cells.Item(2538, "E").Value = 'admin'
$ws.Cells.Item(2538, "F").Value = 'now()'

# row 2539
$ws.Cells.Item(2539, "A").Value = 'ara'
$ws.Cells.Item(2539, "B").Value = 'WEEKLY_SUMMARY_SUBJECT'
$ws.Cells.Item(2539, "C").Value = 'قالب موضوعي لإشعارات الملخص الأسبوعي'
$ws.Cells.Item(2, 4).Copy($ws.Cells.Item(2539, 4)) | Out-Null
$ws.Cells.Item(2539, "E").Value = 'admin'
$ws.Cells.Item(2539, "F").Value = 'now()'

# row 2540
$ws.Cells.Item(2540, "A").Value = 'hin'
$ws.Cells.Item(2540, "B").Value = 'WEEKLY_SUMMARY_SUBJECT'
$ws.Cells.Item(2540, "C").Value = 'साप्ताहिक सारांश अधिसूचनाओं के लिए विषय टेम्पलेट'
$ws.Cells.Item(2, 4).Copy($ws.Cells.Item(2540, 4)) | Out-Null
$ws.Cells.Item(2540, "E").Value = 'admin'
$ws.Cells.Item(2540, "F").Value = 'now()'

# row 2541
$ws.Cells.Item(2541, "A").Value = 'kan'
$ws.Cells.Item(2541, "B").Value = 'WEEKLY_SUMMARY_SUBJECT'
$ws.Cells.Item(2541, "C").Value = 'ವಾರದ ಸಾರಾಂಶ ಅಧಿಸೂಚನೆಗಳಿಗಾಗಿ ವಿಷಯ ಟೆಂಪ್ಲೇಟ್'
$ws.Cells.Item(2, 4).Copy($ws.Cells.Item(2541, 4)) | Out-Null
$ws.Cells.Item(2541, "E").Value = 'admin'
$ws.Cells.Item(2541, "F").Value = 'now()'

# row 2542
$ws.Cells.Item(2542, "A").Value = 'tam'
$ws.Cells.Item(2542, "B").Value = 'WEEKLY_SUMMARY_SUBJECT'
$ws.Cells.Item(2542, "C").Value = 'வாராந்திர சுருக்க அறிவிப்புகளுக்கான தலைப்பு டெம்ப்ளேட்'
$ws.Cells.Item(2, 4).Copy($ws.Cells.Item(2542, 4)) | Out-Null
$ws.Cells.Item(2542, "E").Value = 'admin'
$ws.Cells.Item(2542, "F").Value = 'now()'

# --- View state: scroll position + selection, best effort ---
$excel.ActiveWindow.ScrollRow = 2489
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B2541").Select() | Out-Null

